$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04856815005294025
$ws.Range("C2").Value = 0.00670943022877594

$ws.Range("B3").Value = 0.04041380344049289
$ws.Range("C3").Value = 0.0009189705784143857

$ws.Range("B4").Value = 0.1031532087561053
$ws.Range("C4").Value = 0.01055978047917941

$ws.Range("B5").Value = 0.09081307540638102
$ws.Range("C5").Value = 0.00249664836933218

$ws.Range("B6").Value = 0.09873286292462449
$ws.Range("C6").Value = 0.001227074577470893
